$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2:E9").Value = "ossywy1045687679"

$ws.Range("I2:I9").NumberFormat = "@"
$ws.Range("I2:I9").Value = "2145512363"

$ws.Range("AX2:AX9").NumberFormat = "@"
$ws.Range("AX2:AX9").Value = "0494063210"

$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "9O"

$ws.Range("AU1").Select()
$ws.Range("BD4:BD5").Select()
